# Insert a new weekly price-record row before row 5 (shifts rows 5..58 down
# to 6..59) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'Femacal de La Calera'
$ws.Range("C5").Value = 'Coquimbo'
$ws.Range("D5").Value = 44761
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 100112035
$ws.Range("G5").Value = 'Bruselas (repollito)'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("N5").Value = '$/malla 15 kilos'
$ws.Range("O5").Value = 'Provincia de Quillota'
$ws.Range("P5").Value = 1000
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = 'Hortaliza'
